# Auto-generated edit script applying the Ifrit_Profits diff
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1893.3684
$ws.Range("I19").Value = 3437.75
$ws.Range("J19").Value = 770.1818
$ws.Range("K19").Value = 3437.75
$ws.Range("L19").Value = 770.1818
$ws.Range("M19").Value = -3262.75
$ws.Range("N19").Value = -1120.1818
$ws.Range("H76").Value = 3209.0908
$ws.Range("I76").Value = 3213.9534
$ws.Range("J76").Value = 3000
$ws.Range("K76").Value = 3213.9534
$ws.Range("L76").Value = 3000
$ws.Range("M76").Value = -2898.9534
$ws.Range("N76").Value = -3630
$ws.Range("H79").Value = 3209.0908
$ws.Range("I79").Value = 3213.9534
$ws.Range("J79").Value = 3000
$ws.Range("K79").Value = 3213.9534
$ws.Range("L79").Value = 3000
$ws.Range("M79").Value = -2121.9534
$ws.Range("N79").Value = -5184
$ws.Range("I100").Value = 1836.3636
$ws.Range("K100").Value = 1836.3636
$ws.Range("M100").Value = -1295.3636
$ws.Range("H106").Value = 3320
$ws.Range("I106").Value = 2980
$ws.Range("J106").Value = 4000
$ws.Range("K106").Value = 2980
$ws.Range("L106").Value = 4000
$ws.Range("M106").Value = -2349
$ws.Range("N106").Value = -5262
$ws.Range("H129").Value = 1183.258
$ws.Range("J129").Value = 1339.2307
$ws.Range("L129").Value = 4017.6921
$ws.Range("N129").Value = -14017.6921
$ws.Range("H138").Value = 1174534
$ws.Range("I138").Value = 2453.75
$ws.Range("J138").Value = 2169026.2
$ws.Range("K138").Value = 7361.25
$ws.Range("L138").Value = 6507078.600000001
$ws.Range("M138").Value = -2221.25
$ws.Range("N138").Value = -6517358.600000001
$ws.Range("H141").Value = 2977.2727
$ws.Range("I141").Value = 2977.2727
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 8931.8181
$ws.Range("L141").Value = 0
$ws.Range("M141").Value = -3751.8181
$ws.Range("N141").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1838630.4
$ws.Range("J2").Value = 3676858
$ws.Range("L2").Value = 3676858
$ws.Range("N2").Value = -3677084
$ws.Range("H32").Value = 9635.164000000001
$ws.Range("I32").Value = 5468.12
$ws.Range("J32").Value = 28576.273
$ws.Range("K32").Value = 5468.12
$ws.Range("L32").Value = 28576.273
$ws.Range("M32").Value = -5181.12
$ws.Range("N32").Value = -29150.273
$ws.Range("H61").Value = 2516.8647
$ws.Range("I61").Value = 2253
$ws.Range("J61").Value = 3229.3
$ws.Range("K61").Value = 2253
$ws.Range("L61").Value = 3229.3
$ws.Range("M61").Value = -2041
$ws.Range("N61").Value = -3653.3
$ws.Range("H74").Value = 2635088.8
$ws.Range("I74").Value = 3704247.2
$ws.Range("K74").Value = 3704247.2
$ws.Range("M74").Value = -3703373.2
$ws.Range("H77").Value = 2635088.8
$ws.Range("I77").Value = 3704247.2
$ws.Range("K77").Value = 18521236
$ws.Range("M77").Value = -18516868
$ws.Range("H105").Value = 48000
$ws.Range("J105").Value = 48000
$ws.Range("L105").Value = 48000
$ws.Range("N105").Value = -54988
$ws.Range("H116").Value = 1838630.4
$ws.Range("J116").Value = 3676858
$ws.Range("L116").Value = 3676858
$ws.Range("N116").Value = -3681446
$ws.Range("H136").Value = 2516.8647
$ws.Range("I136").Value = 2253
$ws.Range("J136").Value = 3229.3
$ws.Range("K136").Value = 6759
$ws.Range("L136").Value = 9687.900000000001
$ws.Range("M136").Value = -4209
$ws.Range("N136").Value = -14787.9

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1838630.4
$ws.Range("J3").Value = 3676858
$ws.Range("L3").Value = 3676858
$ws.Range("N3").Value = -3677086
$ws.Range("H98").Value = 44999.5
$ws.Range("J98").Value = 44999.5
$ws.Range("L98").Value = 44999.5
$ws.Range("N98").Value = -50989.5
$ws.Range("H99").Value = 773.3333
$ws.Range("I99").Value = 610
$ws.Range("J99").Value = 1100
$ws.Range("K99").Value = 610
$ws.Range("L99").Value = 1100
$ws.Range("M99").Value = 888
$ws.Range("N99").Value = -4096
$ws.Range("H100").Value = 25000
$ws.Range("J100").Value = 25000
$ws.Range("L100").Value = 25000
$ws.Range("N100").Value = -27164

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H43").Value = 0
$ws.Range("J43").Value = 0
$ws.Range("L43").Value = 0
$ws.Range("N43").ClearContents()
$ws.Range("H101").Value = 0
$ws.Range("J101").Value = 0
$ws.Range("L101").Value = 0
$ws.Range("N101").ClearContents()
$ws.Range("H102").Value = 30541
$ws.Range("J102").Value = 30541
$ws.Range("L102").Value = 30541
$ws.Range("N102").Value = -35409
$ws.Range("H103").Value = 23205.9
$ws.Range("I103").Value = 10834.8
$ws.Range("J103").Value = 35577
$ws.Range("K103").Value = 10834.8
$ws.Range("L103").Value = 35577
$ws.Range("M103").Value = -9662.799999999999
$ws.Range("N103").Value = -37921

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 816.3819999999999
$ws.Range("I68").Value = 609.67346
$ws.Range("J68").Value = 1069.6
$ws.Range("K68").Value = 1829.02038
$ws.Range("L68").Value = 3208.8
$ws.Range("M68").Value = -1018.02038
$ws.Range("N68").Value = -4830.799999999999
$ws.Range("H71").Value = 816.3819999999999
$ws.Range("I71").Value = 609.67346
$ws.Range("J71").Value = 1069.6
$ws.Range("K71").Value = 5487.06114
$ws.Range("L71").Value = 9626.4
$ws.Range("M71").Value = -1431.06114
$ws.Range("N71").Value = -17738.4
$ws.Range("H99").Value = 1142.7142
$ws.Range("I99").Value = 849.8333
$ws.Range("J99").Value = 2900
$ws.Range("K99").Value = 2549.4999
$ws.Range("L99").Value = 8700
$ws.Range("M99").Value = -303.4998999999998
$ws.Range("N99").Value = -13192
$ws.Range("H100").Value = 1980
$ws.Range("J100").Value = 0
$ws.Range("L100").Value = 0
$ws.Range("N100").ClearContents()
$ws.Range("H117").Value = 2413.8462
$ws.Range("I117").Value = 875
$ws.Range("J117").Value = 3097.7778
$ws.Range("K117").Value = 2625
$ws.Range("L117").Value = 9293.3334
$ws.Range("M117").Value = 817
$ws.Range("N117").Value = -16177.3334
$ws.Range("H131").Value = 1820477.6
$ws.Range("J131").Value = 2566958.5
$ws.Range("L131").Value = 7700875.5
$ws.Range("N131").Value = -7710955.5
$ws.Range("H132").Value = 55556612
$ws.Range("I132").Value = 83334520
$ws.Range("J132").Value = 789.8333
$ws.Range("K132").Value = 750010680
$ws.Range("L132").Value = 7108.4997
$ws.Range("M132").Value = -750008150
$ws.Range("N132").Value = -12168.4997

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 862.3333
$ws.Range("I97").Value = 921.7273
$ws.Range("J97").Value = 699
$ws.Range("K97").Value = 921.7273
$ws.Range("L97").Value = 699
$ws.Range("M97").Value = -425.7273
$ws.Range("N97").Value = -1691
$ws.Range("H98").Value = 55000
$ws.Range("J98").Value = 55000
$ws.Range("L98").Value = 55000
$ws.Range("N98").Value = -60990
$ws.Range("H99").Value = 1768.3
$ws.Range("I99").Value = 1768.3
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 1768.3
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = 477.7
$ws.Range("N99").ClearContents()
$ws.Range("H132").Value = 25643504
$ws.Range("I132").Value = 55558276
$ws.Range("J132").Value = 2271
$ws.Range("K132").Value = 166674828
$ws.Range("L132").Value = 6813
$ws.Range("M132").Value = -166672298
$ws.Range("N132").Value = -11873

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H101").Value = 19362
$ws.Range("J101").Value = 19362
$ws.Range("L101").Value = 19362
$ws.Range("N101").Value = -25852
$ws.Range("H102").Value = 48000
$ws.Range("J102").Value = 48000
$ws.Range("L102").Value = 48000
$ws.Range("N102").Value = -54490

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H98").Value = 19084.285
$ws.Range("J98").Value = 19084.285
$ws.Range("L98").Value = 19084.285
$ws.Range("N98").Value = -25074.285
